$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 25392.95
$ws.Range("J12").Value = 63062.375
$ws.Range("L12").Value = 63062.375
$ws.Range("N12").Value = -63402.375
$ws.Range("H64").Value = 74034.5
$ws.Range("I64").Value = 127400
$ws.Range("J64").Value = 2880.5
$ws.Range("K64").Value = 127400
$ws.Range("L64").Value = 2880.5
$ws.Range("M64").Value = -127152
$ws.Range("N64").Value = -3376.5
$ws.Range("H67").Value = 74034.5
$ws.Range("I67").Value = 127400
$ws.Range("J67").Value = 2880.5
$ws.Range("K67").Value = 127400
$ws.Range("L67").Value = 2880.5
$ws.Range("M67").Value = -126542
$ws.Range("N67").Value = -4596.5
$ws.Range("H70").Value = 2063.6365
$ws.Range("I70").Value = 2322.2222
$ws.Range("K70").Value = 6966.6666
$ws.Range("M70").Value = -6696.6666
$ws.Range("H73").Value = 2063.6365
$ws.Range("I73").Value = 2322.2222
$ws.Range("K73").Value = 6966.6666
$ws.Range("M73").Value = -6030.6666
$ws.Range("H103").Value = 1085.9259
$ws.Range("I103").Value = 2950
$ws.Range("J103").Value = 761.73914
$ws.Range("K103").Value = 8850
$ws.Range("L103").Value = 2285.21742
$ws.Range("M103").Value = -8264
$ws.Range("N103").Value = -3457.21742
$ws.Range("H137").Value = 1336.5758
$ws.Range("I137").Value = 1222.0938
$ws.Range("K137").Value = 3666.2814
$ws.Range("M137").Value = -1116.2814
$ws.Range("H138").Value = 1307.4529
$ws.Range("I138").Value = 1047.7609
$ws.Range("J138").Value = 3014
$ws.Range("K138").Value = 3143.2827
$ws.Range("L138").Value = 9042
$ws.Range("M138").Value = 1996.7173
$ws.Range("N138").Value = -19322

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 43891.25
$ws.Range("I102").Value = 78988.38
$ws.Range("J102").Value = 2412.818
$ws.Range("K102").Value = 78988.38
$ws.Range("L102").Value = 2412.818
$ws.Range("M102").Value = -77366.38
$ws.Range("N102").Value = -5656.818
$ws.Range("I110").Value = 71578990
$ws.Range("K110").Value = 71578990
$ws.Range("M110").Value = -71576945

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 598.0833
$ws.Range("I94").Value = 489.64706
$ws.Range("K94").Value = 489.64706
$ws.Range("M94").Value = -38.64706000000001
$ws.Range("H99").Value = 1955.68
$ws.Range("I99").Value = 1631.8334
$ws.Range("J99").Value = 2057.9473
$ws.Range("K99").Value = 1631.8334
$ws.Range("L99").Value = 2057.9473
$ws.Range("M99").Value = -133.8334
$ws.Range("N99").Value = -5053.9473
$ws.Range("H107").Value = 62553940
$ws.Range("J107").Value = 404.66666
$ws.Range("L107").Value = 404.66666
$ws.Range("N107").Value = -4244.66666

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 20240.875
$ws.Range("J99").Value = 34014
$ws.Range("L99").Value = 34014
$ws.Range("N99").Value = -37010
$ws.Range("H126").Value = 20240.875
$ws.Range("J126").Value = 34014
$ws.Range("L126").Value = 102042
$ws.Range("N126").Value = -106982

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 43.545456
$ws.Range("I12").Value = 39
$ws.Range("J12").Value = 46.142857
$ws.Range("K12").Value = 117
$ws.Range("L12").Value = 138.428571
$ws.Range("M12").Value = 56
$ws.Range("N12").Value = -484.428571
$ws.Range("H23").Value = 720.8570999999999
$ws.Range("I23").Value = 500
$ws.Range("J23").Value = 737.8461
$ws.Range("K23").Value = 1500
$ws.Range("L23").Value = 2213.5383
$ws.Range("M23").Value = -1265
$ws.Range("N23").Value = -2683.5383
$ws.Range("H33").Value = 2333.4443
$ws.Range("J33").Value = 2333.4443
$ws.Range("L33").Value = 14000.6658
$ws.Range("N33").Value = -14566.6658
$ws.Range("H38").Value = 53.6
$ws.Range("I38").Value = 15
$ws.Range("J38").Value = 63.25
$ws.Range("K38").Value = 45
$ws.Range("L38").Value = 189.75
$ws.Range("M38").Value = 302
$ws.Range("N38").Value = -883.75
$ws.Range("H44").Value = 721
$ws.Range("I44").Value = 243.71428
$ws.Range("J44").Value = 1834.6666
$ws.Range("K44").Value = 731.14284
$ws.Range("L44").Value = 5503.9998
$ws.Range("M44").Value = -333.14284
$ws.Range("N44").Value = -6299.9998
$ws.Range("H105").Value = 7453.8
$ws.Range("J105").Value = 7453.8
$ws.Range("L105").Value = 22361.4
$ws.Range("N105").Value = -27603.4

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2741.1738
$ws.Range("I102").Value = 2500.6155
$ws.Range("J102").Value = 3053.9
$ws.Range("K102").Value = 2500.6155
$ws.Range("L102").Value = 3053.9
$ws.Range("M102").Value = -878.6154999999999
$ws.Range("N102").Value = -6297.9

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3705.0557
$ws.Range("I7").Value = 2515
$ws.Range("K7").Value = 2515
$ws.Range("M7").Value = -2403
$ws.Range("H61").Value = 1882.091
$ws.Range("I61").Value = 1890.3
$ws.Range("J61").Value = 1800
$ws.Range("K61").Value = 1890.3
$ws.Range("L61").Value = 1800
$ws.Range("M61").Value = -1688.3
$ws.Range("N61").Value = -2204
$ws.Range("H68").Value = 3694.6155
$ws.Range("I68").Value = 1725
$ws.Range("J68").Value = 4570
$ws.Range("K68").Value = 1725
$ws.Range("L68").Value = 4570
$ws.Range("M68").Value = -976
$ws.Range("N68").Value = -6068
$ws.Range("H71").Value = 3694.6155
$ws.Range("I71").Value = 1725
$ws.Range("J71").Value = 4570
$ws.Range("K71").Value = 8625
$ws.Range("L71").Value = 22850
$ws.Range("M71").Value = -4881
$ws.Range("N71").Value = -30338
$ws.Range("H82").Value = 1989.4445
$ws.Range("I82").Value = 987.5
$ws.Range("J82").Value = 2791
$ws.Range("K82").Value = 987.5
$ws.Range("L82").Value = 2791
$ws.Range("M82").Value = -626.5
$ws.Range("N82").Value = -3513
$ws.Range("H85").Value = 1989.4445
$ws.Range("I85").Value = 987.5
$ws.Range("J85").Value = 2791
$ws.Range("K85").Value = 987.5
$ws.Range("L85").Value = 2791
$ws.Range("M85").Value = 260.5
$ws.Range("N85").Value = -5287
$ws.Range("H93").Value = 1605.9333
$ws.Range("I93").Value = 1782.1666
$ws.Range("J93").Value = 901
$ws.Range("K93").Value = 1782.1666
$ws.Range("L93").Value = 901
$ws.Range("M93").Value = -534.1666
$ws.Range("N93").Value = -3397
$ws.Range("H100").Value = 2810.889
$ws.Range("I100").Value = 1779.8
$ws.Range("J100").Value = 4099.75
$ws.Range("K100").Value = 1779.8
$ws.Range("L100").Value = 4099.75
$ws.Range("M100").Value = -1238.8
$ws.Range("N100").Value = -5181.75
$ws.Range("H113").Value = 1882.091
$ws.Range("I113").Value = 1890.3
$ws.Range("J113").Value = 1800
$ws.Range("K113").Value = 1890.3
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = 279.7
$ws.Range("N113").Value = -6140
$ws.Range("H122").Value = 2491.4546
$ws.Range("I122").Value = 2491.4546
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7474.3638
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5024.3638
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 3705.0557
$ws.Range("I126").Value = 2515
$ws.Range("K126").Value = 7545
$ws.Range("M126").Value = -5075

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 62501948
$ws.Range("I96").Value = 166669020
$ws.Range("J96").Value = 1700.4
$ws.Range("K96").Value = 166669020
$ws.Range("L96").Value = 1700.4
$ws.Range("M96").Value = -166667647
$ws.Range("N96").Value = -4446.4
$ws.Range("H122").Value = 3249.9167
$ws.Range("I122").Value = 6152
$ws.Range("J122").Value = 2669.5
$ws.Range("K122").Value = 18456
$ws.Range("L122").Value = 8008.5
$ws.Range("M122").Value = -16006
$ws.Range("N122").Value = -12908.5
$ws.Range("H126").Value = 1995.1
$ws.Range("I126").Value = 1923
$ws.Range("K126").Value = 5769
$ws.Range("M126").Value = -3299
$ws.Range("H136").Value = 629.64813
$ws.Range("I136").Value = 362.15555
$ws.Range("J136").Value = 1967.1111
$ws.Range("K136").Value = 1086.46665
$ws.Range("L136").Value = 5901.3333
$ws.Range("M136").Value = 1463.53335
$ws.Range("N136").Value = -11001.3333
